$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new blank rows before row 165, shifting existing rows 165:184 down to 168:187
$ws.Range("A165:T167").Insert()

# Row 165: Dina / Especial
$ws.Range("A165").Value = 3
$ws.Range("B165").Value = "Femacal de La Calera"
$ws.Range("C165").Value = "Coquimbo"
$ws.Range("D165").Value = 44918
$ws.Range("E165").Value = 5
$ws.Range("F165").Value = "Fruta"
$ws.Range("G165").Value = 100103
$ws.Range("H165").Value = "Frutos de hueso (carozo)"
$ws.Range("I165").Value = 100103003
$ws.Range("J165").Value = "Damasco"
$ws.Range("K165").Value = "Dina"
$ws.Range("L165").Value = "Especial"
$ws.Range("M165").Value = 56
$ws.Range("N165").Value = 18000
$ws.Range("O165").Value = 18000
$ws.Range("P165").Value = 18000
$ws.Range("Q165").Value = "$/caja 15 kilos"
$ws.Range("R165").Value = "Región de O'Higgins"
$ws.Range("S165").Value = 1200
$ws.Range("T165").Value = 15

# Row 166: Dina / Primera
$ws.Range("A166").Value = 3
$ws.Range("B166").Value = "Femacal de La Calera"
$ws.Range("C166").Value = "Coquimbo"
$ws.Range("D166").Value = 44918
$ws.Range("E166").Value = 5
$ws.Range("F166").Value = "Fruta"
$ws.Range("G166").Value = 100103
$ws.Range("H166").Value = "Frutos de hueso (carozo)"
$ws.Range("I166").Value = 100103003
$ws.Range("J166").Value = "Damasco"
$ws.Range("K166").Value = "Dina"
$ws.Range("L166").Value = "Primera"
$ws.Range("M166").Value = 60
$ws.Range("N166").Value = 15000
$ws.Range("O166").Value = 15000
$ws.Range("P166").Value = 15000
$ws.Range("Q166").Value = "$/caja 15 kilos"
$ws.Range("R166").Value = "Región de O'Higgins"
$ws.Range("S166").Value = 1000
$ws.Range("T166").Value = 15

# Row 167: Dina / Segunda
$ws.Range("A167").Value = 3
$ws.Range("B167").Value = "Femacal de La Calera"
$ws.Range("C167").Value = "Coquimbo"
$ws.Range("D167").Value = 44918
$ws.Range("E167").Value = 5
$ws.Range("F167").Value = "Fruta"
$ws.Range("G167").Value = 100103
$ws.Range("H167").Value = "Frutos de hueso (carozo)"
$ws.Range("I167").Value = 100103003
$ws.Range("J167").Value = "Damasco"
$ws.Range("K167").Value = "Dina"
$ws.Range("L167").Value = "Segunda"
$ws.Range("M167").Value = 60
$ws.Range("N167").Value = 11000
$ws.Range("O167").Value = 11000
$ws.Range("P167").Value = 11000
$ws.Range("Q167").Value = "$/caja 15 kilos"
$ws.Range("R167").Value = "Región de O'Higgins"
$ws.Range("S167").Value = 733
$ws.Range("T167").Value = 15
